# kick PARTNER from project group for AUTHORs and ADMINs
#
# 1. Insert two new columns ("speciality", "group") before the existing
#    "phonenum" column (K), shifting phonenum/email/status/projectMember
#    two columns to the right (K->M, L->N, M->O, N->P).
# 2. Fill in the new header cells and the new row-2 data.
# 3. Update row 2's "aboutme" text, and its (now shifted) status/
#    projectMember values.
# 4. Swap the usernames of rows 3 and 4, and drop the stray lastname
#    value that used to tag along with "riba17".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the two new columns at K:L (pushes K.. right by 2) ---
$ws.Range("K1:L1").EntireColumn.Insert()

# --- headers for the new columns ---
$ws.Range("K1").Value = "speciality"
$ws.Range("L1").Value = "group"

# --- row 2 updates ---
$ws.Range("F2").Value = "Программист"
$ws.Range("K2").Value = "Инжиниринг предприятий и информационных систем"
$ws.Range("L2").Value = "ПИ03у"
$ws.Range("O2").Value = "{RESIDENT,ADMIN,EVENT_MANAGER}"
$ws.Range("P2").Value = "TimeTrace"

# --- rows 3 & 4: swap usernames, drop the leftover lastname on row 3 ---
$ws.Range("A3").Value = "whereistheexit"
$ws.Range("A4").Value = "riba17"
$ws.Range("B3").ClearContents()
